$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "Ejercicio 4"

$ws.Range("A17").Value = "Nivel"

$ws.Range("A18").Value = 1
$ws.Range("A19").Value = 2
$ws.Range("A20").Value = 3
$ws.Range("A21").Value = 4

$ws.Range("B18").Value = 1
$ws.Range("B19").Value = 5
$ws.Range("B20").Value = 17
$ws.Range("B21").Value = 53

$ws.Range("C19").Value = "f(n-1) +4"
$ws.Range("C20").Value = "f(n-1) + 12"
$ws.Range("C21").Value = "f(n-1) + 36"

$ws.Range("D19").Value = "f(n-1) + 4 * 1"
$ws.Range("D20").Value = "f(n-1) + 4  * 3"
$ws.Range("D21").Value = "f(n-1) + 4 * 9"

$ws.Range("E19").Value = "f(n-1) + 4 * 3^0"
$ws.Range("E19").Font.Underline = $true
$ws.Range("E20").Value = "f(n-1) + 4 * 3^1"
$ws.Range("E21").Value = "f(n-1) + 4 * 3^2"

$ws.Range("F19").Value = "f(n-1) + 4 * 3^(n-2)"
$ws.Range("F20").Value = "f(n-1) + 4 * 3^(n-2)"
$ws.Range("F21").Value = "f(n-1) + 4 * 3^(n-2)"

$ws.Columns.Item(5).ColumnWidth = 17

$ws.Range("F21").Select() | Out-Null
